$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Uts2"
$ws.Range("C2").Value = "Uts2r"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.176339
$ws.Range("H2").Value = 0.529017
$ws.Range("I2").Value = 0.2427660118313651
$ws.Range("J2").Value = 0.2427660118313651
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.024594
$ws.Range("N2").Value = 0.073782
$ws.Range("O2").Value = 0.07315455339552658
$ws.Range("P2").Value = 0.07315455339552657
$ws.Range("Q2").Value = 0.004336881366
$ws.Range("R2").Value = 0.039031932294
$ws.Range("S2").Value = 0.01775943917513664
$ws.Range("T2").Value = 0.01775943917513664

$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Uts2"
$ws.Range("C3").Value = "Uts2r"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.176339
$ws.Range("H3").Value = 0.529017
$ws.Range("I3").Value = 0.2427660118313651
$ws.Range("J3").Value = 0.2427660118313651
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.2813096666666666
$ws.Range("N3").Value = 0.8439289999999999
$ws.Range("O3").Value = 0.8367521765814607
$ws.Range("P3").Value = 0.8367521765814607
$ws.Range("Q3").Value = 0.04960586531033333
$ws.Range("R3").Value = 0.4464527877929999
$ws.Range("S3").Value = 0.2031349887998954
$ws.Range("T3").Value = 0.2031349887998954

$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Uts2"
$ws.Range("C4").Value = "Uts2r"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.176339
$ws.Range("H4").Value = 0.529017
$ws.Range("I4").Value = 0.2427660118313651
$ws.Range("J4").Value = 0.2427660118313651
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.01219966666666667
$ws.Range("N4").Value = 0.036599
$ws.Range("O4").Value = 0.0362877598834794
$ws.Range("P4").Value = 0.0362877598834794
$ws.Range("Q4").Value = 0.002151277020333333
$ws.Range("R4").Value = 0.019361493183
$ws.Range("S4").Value = 0.008809434745206497
$ws.Range("T4").Value = 0.008809434745206497

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Uts2"
$ws.Range("C5").Value = "Uts2r"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.176339
$ws.Range("H5").Value = 0.529017
$ws.Range("I5").Value = 0.2427660118313651
$ws.Range("J5").Value = 0.2427660118313651
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.018089
$ws.Range("N5").Value = 0.054267
$ws.Range("O5").Value = 0.05380551013953323
$ws.Range("P5").Value = 0.05380551013953323
$ws.Range("Q5").Value = 0.003189796171
$ws.Range("R5").Value = 0.028708165539
$ws.Range("S5").Value = 0.01306214911112656
$ws.Range("T5").Value = 0.01306214911112656

$ws.Range("A6").Value = "Inflammatory-Mac"
$ws.Range("B6").Value = "Uts2"
$ws.Range("C6").Value = "Uts2r"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.2935233333333334
$ws.Range("H6").Value = 0.8805700000000001
$ws.Range("I6").Value = 0.4040937569838875
$ws.Range("J6").Value = 0.4040937569838876
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.024594
$ws.Range("N6").Value = 0.073782
$ws.Range("O6").Value = 0.07315455339552658
$ws.Range("P6").Value = 0.07315455339552657
$ws.Range("Q6").Value = 0.007218912860000001
$ws.Range("R6").Value = 0.06497021574
$ws.Range("S6").Value = 0.02956129832207674
$ws.Range("T6").Value = 0.02956129832207674

$ws.Range("A7").Value = "Inflammatory-Mac"
$ws.Range("B7").Value = "Uts2"
$ws.Range("C7").Value = "Uts2r"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.2935233333333334
$ws.Range("H7").Value = 0.8805700000000001
$ws.Range("I7").Value = 0.4040937569838875
$ws.Range("J7").Value = 0.4040937569838876
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.2813096666666666
$ws.Range("N7").Value = 0.8439289999999999
$ws.Range("O7").Value = 0.8367521765814607
$ws.Range("P7").Value = 0.8367521765814607
$ws.Range("Q7").Value = 0.08257095105888888
$ws.Range("R7").Value = 0.74313855953
$ws.Range("S7").Value = 0.3381263306992477
$ws.Range("T7").Value = 0.3381263306992477

$ws.Range("A8").Value = "Inflammatory-Mac"
$ws.Range("B8").Value = "Uts2"
$ws.Range("C8").Value = "Uts2r"
$ws.Range("D8").Value = "MuSCs"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.2935233333333334
$ws.Range("H8").Value = 0.8805700000000001
$ws.Range("I8").Value = 0.4040937569838875
$ws.Range("J8").Value = 0.4040937569838876
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.01219966666666667
$ws.Range("N8").Value = 0.036599
$ws.Range("O8").Value = 0.0362877598834794
$ws.Range("P8").Value = 0.0362877598834794
$ws.Range("Q8").Value = 0.003580886825555555
$ws.Range("R8").Value = 0.03222798143
$ws.Range("S8").Value = 0.01466365722384439
$ws.Range("T8").Value = 0.01466365722384439

$ws.Range("A9").Value = "Inflammatory-Mac"
$ws.Range("B9").Value = "Uts2"
$ws.Range("C9").Value = "Uts2r"
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.2935233333333334
$ws.Range("H9").Value = 0.8805700000000001
$ws.Range("I9").Value = 0.4040937569838875
$ws.Range("J9").Value = 0.4040937569838876
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.018089
$ws.Range("N9").Value = 0.054267
$ws.Range("O9").Value = 0.05380551013953323
$ws.Range("P9").Value = 0.05380551013953323
$ws.Range("Q9").Value = 0.005309543576666667
$ws.Range("R9").Value = 0.04778589219000001
$ws.Range("S9").Value = 0.02174247073871863
$ws.Range("T9").Value = 0.02174247073871864

$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Uts2"
$ws.Range("C10").Value = "Uts2r"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.1622493333333333
$ws.Range("H10").Value = 0.486748
$ws.Range("I10").Value = 0.2233687588997959
$ws.Range("J10").Value = 0.2233687588997959
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.024594
$ws.Range("N10").Value = 0.073782
$ws.Range("O10").Value = 0.07315455339552658
$ws.Range("P10").Value = 0.07315455339552657
$ws.Range("Q10").Value = 0.003990360104
$ws.Range("R10").Value = 0.035913240936
$ws.Range("S10").Value = 0.01634044179982762
$ws.Range("T10").Value = 0.01634044179982762

$ws.Range("A11").Value = "MuSCs"
$ws.Range("B11").Value = "Uts2"
$ws.Range("C11").Value = "Uts2r"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 0.6666666666666666
$ws.Range("G11").Value = 0.1622493333333333
$ws.Range("H11").Value = 0.486748
$ws.Range("I11").Value = 0.2233687588997959
$ws.Range("J11").Value = 0.2233687588997959
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.2813096666666666
$ws.Range("N11").Value = 0.8439289999999999
$ws.Range("O11").Value = 0.8367521765814607
$ws.Range("P11").Value = 0.8367521765814607
$ws.Range("Q11").Value = 0.04564230587688888
$ws.Range("R11").Value = 0.4107807528919999
$ws.Range("S11").Value = 0.1869042951897037
$ws.Range("T11").Value = 0.1869042951897038

$ws.Range("A12").Value = "MuSCs"
$ws.Range("B12").Value = "Uts2"
$ws.Range("C12").Value = "Uts2r"
$ws.Range("D12").Value = "MuSCs"
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 0.6666666666666666
$ws.Range("G12").Value = 0.1622493333333333
$ws.Range("H12").Value = 0.486748
$ws.Range("I12").Value = 0.2233687588997959
$ws.Range("J12").Value = 0.2233687588997959
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 0.6666666666666666
$ws.Range("M12").Value = 0.01219966666666667
$ws.Range("N12").Value = 0.036599
$ws.Range("O12").Value = 0.0362877598834794
$ws.Range("P12").Value = 0.0362877598834794
$ws.Range("Q12").Value = 0.001979387783555555
$ws.Range("R12").Value = 0.017814490052
$ws.Range("S12").Value = 0.008105551888426593
$ws.Range("T12").Value = 0.008105551888426595

$ws.Range("A13").Value = "MuSCs"
$ws.Range("B13").Value = "Uts2"
$ws.Range("C13").Value = "Uts2r"
$ws.Range("D13").Value = "Resolving-Mac"
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 0.6666666666666666
$ws.Range("G13").Value = 0.1622493333333333
$ws.Range("H13").Value = 0.486748
$ws.Range("I13").Value = 0.2233687588997959
$ws.Range("J13").Value = 0.2233687588997959
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 0.3333333333333333
$ws.Range("M13").Value = 0.018089
$ws.Range("N13").Value = 0.054267
$ws.Range("O13").Value = 0.05380551013953323
$ws.Range("P13").Value = 0.05380551013953323
$ws.Range("Q13").Value = 0.002934928190666667
$ws.Range("R13").Value = 0.026414353716
$ws.Range("S13").Value = 0.01201847002183792
$ws.Range("T13").Value = 0.01201847002183792

$ws.Range("A14").Value = "Resolving-Mac"
$ws.Range("B14").Value = "Uts2"
$ws.Range("C14").Value = "Uts2r"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.09426266666666666
$ws.Range("H14").Value = 0.282788
$ws.Range("I14").Value = 0.1297714722849513
$ws.Range("J14").Value = 0.1297714722849513
$ws.Range("K14").Value = 1
$ws.Range("L14").Value = 0.3333333333333333
$ws.Range("M14").Value = 0.024594
$ws.Range("N14").Value = 0.073782
$ws.Range("O14").Value = 0.07315455339552658
$ws.Range("P14").Value = 0.07315455339552657
$ws.Range("Q14").Value = 0.002318296024
$ws.Range("R14").Value = 0.020864664216
$ws.Range("S14").Value = 0.00949337409848557
$ws.Range("T14").Value = 0.009493374098485568

$ws.Range("A15").Value = "Resolving-Mac"
$ws.Range("B15").Value = "Uts2"
$ws.Range("C15").Value = "Uts2r"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.09426266666666666
$ws.Range("H15").Value = 0.282788
$ws.Range("I15").Value = 0.1297714722849513
$ws.Range("J15").Value = 0.1297714722849513
$ws.Range("K15").Value = 2
$ws.Range("L15").Value = 0.6666666666666666
$ws.Range("M15").Value = 0.2813096666666666
$ws.Range("N15").Value = 0.8439289999999999
$ws.Range("O15").Value = 0.8367521765814607
$ws.Range("P15").Value = 0.8367521765814607
$ws.Range("Q15").Value = 0.02651699933911111
$ws.Range("R15").Value = 0.238652994052
$ws.Range("S15").Value = 0.1085865618926137
$ws.Range("T15").Value = 0.1085865618926137

$ws.Range("A16").Value = "Resolving-Mac"
$ws.Range("B16").Value = "Uts2"
$ws.Range("C16").Value = "Uts2r"
$ws.Range("D16").Value = "MuSCs"
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.09426266666666666
$ws.Range("H16").Value = 0.282788
$ws.Range("I16").Value = 0.1297714722849513
$ws.Range("J16").Value = 0.1297714722849513
$ws.Range("K16").Value = 2
$ws.Range("L16").Value = 0.6666666666666666
$ws.Range("M16").Value = 0.01219966666666667
$ws.Range("N16").Value = 0.036599
$ws.Range("O16").Value = 0.0362877598834794
$ws.Range("P16").Value = 0.0362877598834794
$ws.Range("Q16").Value = 0.001149973112444444
$ws.Range("R16").Value = 0.010349758012
$ws.Range("S16").Value = 0.004709116026001914
$ws.Range("T16").Value = 0.004709116026001914

$ws.Range("A17").Value = "Resolving-Mac"
$ws.Range("B17").Value = "Uts2"
$ws.Range("C17").Value = "Uts2r"
$ws.Range("D17").Value = "Resolving-Mac"
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.09426266666666666
$ws.Range("H17").Value = 0.282788
$ws.Range("I17").Value = 0.1297714722849513
$ws.Range("J17").Value = 0.1297714722849513
$ws.Range("K17").Value = 1
$ws.Range("L17").Value = 0.3333333333333333
$ws.Range("M17").Value = 0.018089
$ws.Range("N17").Value = 0.054267
$ws.Range("O17").Value = 0.05380551013953323
$ws.Range("P17").Value = 0.05380551013953323
$ws.Range("Q17").Value = 0.001705117377333333
$ws.Range("R17").Value = 0.015346056396
$ws.Range("S17").Value = 0.006982420267850103
$ws.Range("T17").Value = 0.006982420267850103
